$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('I6').Value = 'ba'
$ws.Range('J6').Value = 'Appreciation'
$ws.Range('I11').Value = 'sv'
$ws.Range('J11').Value = 'Statement-opinion'
$ws.Range('I16').Value = 'ba'
$ws.Range('J16').Value = 'Appreciation'
$ws.Range('I20').Value = 'sv'
$ws.Range('J20').Value = 'Statement-opinion'
$ws.Range('I37').Value = 'sd'
$ws.Range('J37').Value = 'Statement-non-opinion'
$ws.Range('I59').Value = 'sd'
$ws.Range('J59').Value = 'Statement-non-opinion'
$ws.Range('I82').Value = '%'
$ws.Range('J82').Value = 'Uninterpretable'
$ws.Range('I92').Value = 'sv'
$ws.Range('J92').Value = 'Statement-opinion'
$ws.Range('I100').Value = 'aa'
$ws.Range('J100').Value = 'Agree/Accept'
$ws.Range('I101').Value = 'sv'
$ws.Range('J101').Value = 'Statement-opinion'
$ws.Range('I103').Value = 'sd'
$ws.Range('J103').Value = 'Statement-non-opinion'
$ws.Range('I108').Value = 'sd'
$ws.Range('J108').Value = 'Statement-non-opinion'
$ws.Range('I111').Value = 'sv'
$ws.Range('J111').Value = 'Statement-opinion'
$ws.Range('I120').Value = '%'
$ws.Range('J120').Value = 'Uninterpretable'
$ws.Range('I124').Value = 'sd'
$ws.Range('J124').Value = 'Statement-non-opinion'
$ws.Range('I136').Value = 'sd'
$ws.Range('J136').Value = 'Statement-non-opinion'
$ws.Range('I153').Value = 'sd'
$ws.Range('J153').Value = 'Statement-non-opinion'
$ws.Range('I159').Value = 'sd'
$ws.Range('J159').Value = 'Statement-non-opinion'
$ws.Range('I178').Value = 'sv'
$ws.Range('J178').Value = 'Statement-opinion'
$ws.Range('I179').Value = 'sd'
$ws.Range('J179').Value = 'Statement-non-opinion'
$ws.Range('I185').Value = 'sd'
$ws.Range('J185').Value = 'Statement-non-opinion'
$ws.Range('I194').Value = 'ba'
$ws.Range('J194').Value = 'Appreciation'
$ws.Range('I200').Value = 'aa'
$ws.Range('J200').Value = 'Agree/Accept'
$ws.Range('I201').Value = 'ba'
$ws.Range('J201').Value = 'Appreciation'
$ws.Range('I217').Value = 'b'
$ws.Range('J217').Value = 'Acknowledge (Backchannel)'
$ws.Range('I219').Value = 'sd'
$ws.Range('J219').Value = 'Statement-non-opinion'
$ws.Range('I226').Value = 'aa'
$ws.Range('J226').Value = 'Agree/Accept'
$ws.Range('I241').Value = 'aa'
$ws.Range('J241').Value = 'Agree/Accept'
$ws.Range('I242').Value = 'sv'
$ws.Range('J242').Value = 'Statement-opinion'
$ws.Range('I265').Value = '%'
$ws.Range('J265').Value = 'Uninterpretable'
$ws.Range('I271').Value = 'b'
$ws.Range('J271').Value = 'Acknowledge (Backchannel)'
$ws.Range('I288').Value = 'sv'
$ws.Range('J288').Value = 'Statement-opinion'
$ws.Range('I301').Value = 'aa'
$ws.Range('J301').Value = 'Agree/Accept'
$ws.Range('I303').Value = '%'
$ws.Range('J303').Value = 'Uninterpretable'
$ws.Range('I335').Value = 'b'
$ws.Range('J335').Value = 'Acknowledge (Backchannel)'
$ws.Range('I337').Value = 'sd'
$ws.Range('J337').Value = 'Statement-non-opinion'
$ws.Range('I347').Value = 'sd'
$ws.Range('J347').Value = 'Statement-non-opinion'
$ws.Range('I353').Value = 'sd'
$ws.Range('J353').Value = 'Statement-non-opinion'
$ws.Range('I354').Value = 'qy'
$ws.Range('J354').Value = 'Yes-No-Question'
$ws.Range('I355').Value = 'b'
$ws.Range('J355').Value = 'Acknowledge (Backchannel)'
$ws.Range('I356').Value = 'sv'
$ws.Range('J356').Value = 'Statement-opinion'
$ws.Range('I364').Value = 'sd'
$ws.Range('J364').Value = 'Statement-non-opinion'
$ws.Range('I381').Value = 'sd'
$ws.Range('J381').Value = 'Statement-non-opinion'
$ws.Range('I383').Value = 'qy'
$ws.Range('J383').Value = 'Yes-No-Question'
$ws.Range('I388').Value = '%'
$ws.Range('J388').Value = 'Uninterpretable'
$ws.Range('I389').Value = 'sv'
$ws.Range('J389').Value = 'Statement-opinion'
$ws.Range('I402').Value = 'sd'
$ws.Range('J402').Value = 'Statement-non-opinion'
$ws.Range('I404').Value = 'sd'
$ws.Range('J404').Value = 'Statement-non-opinion'
$ws.Range('I407').Value = 'sd'
$ws.Range('J407').Value = 'Statement-non-opinion'
$ws.Range('I409').Value = 'b'
$ws.Range('J409').Value = 'Acknowledge (Backchannel)'
$ws.Range('I413').Value = 'sv'
$ws.Range('J413').Value = 'Statement-opinion'
$ws.Range('I422').Value = '%'
$ws.Range('J422').Value = 'Uninterpretable'
$ws.Range('I423').Value = '%'
$ws.Range('J423').Value = 'Uninterpretable'
$ws.Range('I432').Value = 'sv'
$ws.Range('J432').Value = 'Statement-opinion'
$ws.Range('I434').Value = 'ba'
$ws.Range('J434').Value = 'Appreciation'
$ws.Range('I436').Value = 'sv'
$ws.Range('J436').Value = 'Statement-opinion'
$ws.Range('I443').Value = '%'
$ws.Range('J443').Value = 'Uninterpretable'
$ws.Range('I446').Value = 'sd'
$ws.Range('J446').Value = 'Statement-non-opinion'
$ws.Range('I454').Value = 'sd'
$ws.Range('J454').Value = 'Statement-non-opinion'
$ws.Range('I455').Value = '%'
$ws.Range('J455').Value = 'Uninterpretable'
$ws.Range('I464').Value = '%'
$ws.Range('J464').Value = 'Uninterpretable'
$ws.Range('I467').Value = 'sv'
$ws.Range('J467').Value = 'Statement-opinion'
$ws.Range('I473').Value = 'sd'
$ws.Range('J473').Value = 'Statement-non-opinion'
$ws.Range('I482').Value = 'sv'
$ws.Range('J482').Value = 'Statement-opinion'
$ws.Range('I485').Value = 'sd'
$ws.Range('J485').Value = 'Statement-non-opinion'
$ws.Range('I491').Value = 'ba'
$ws.Range('J491').Value = 'Appreciation'
$ws.Range('I494').Value = '%'
$ws.Range('J494').Value = 'Uninterpretable'
$ws.Range('I515').Value = 'b'
$ws.Range('J515').Value = 'Acknowledge (Backchannel)'
$ws.Range('I523').Value = 'ba'
$ws.Range('J523').Value = 'Appreciation'
$ws.Range('I543').Value = 'sv'
$ws.Range('J543').Value = 'Statement-opinion'
$ws.Range('I547').Value = 'b'
$ws.Range('J547').Value = 'Acknowledge (Backchannel)'
$ws.Range('I551').Value = 'b'
$ws.Range('J551').Value = 'Acknowledge (Backchannel)'
$ws.Range('I557').Value = 'aa'
$ws.Range('J557').Value = 'Agree/Accept'
$ws.Range('I558').Value = '%'
$ws.Range('J558').Value = 'Uninterpretable'
$ws.Range('I560').Value = 'sd'
$ws.Range('J560').Value = 'Statement-non-opinion'
$ws.Range('I571').Value = 'sd'
$ws.Range('J571').Value = 'Statement-non-opinion'
$ws.Range('I577').Value = 'b'
$ws.Range('J577').Value = 'Acknowledge (Backchannel)'
$ws.Range('I579').Value = 'sv'
$ws.Range('J579').Value = 'Statement-opinion'
$ws.Range('I586').Value = 'b'
$ws.Range('J586').Value = 'Acknowledge (Backchannel)'
